$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invoice number
$ws.Range("I3").Value = "A21216407"

# Vendor/store block
$ws.Range("F6").Value = 3
$ws.Range("H6").Value = "Office Depot"

# Line item 1 (row 9)
$ws.Range("A9").Value = "21101"
$ws.Range("B9").Value = "Caja de papel bond"
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = 350
$ws.Range("J9").Value = "Ninguna"

# Line item 2 (row 10) - cleared out entirely
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("J10").Value = ""
